# Applies the "specific aims editing in progress" revision to the F32
# specific aims document using Find/Replace over the Word object model.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# 1. "Deep neural networks have demonstrated the ability to " -> add " (DNN)"
Replace-Text "Deep neural networks have demonstrated the ability to " "Deep neural networks (DNN) have demonstrated the ability to "

# 2. "To most effectively utilize these networks to aid our understanding of splicing " ->
#    "For DNNs to be maximally effective in aiding our understanding of splicing, "
Replace-Text "To most effectively utilize these networks to aid our understanding of splicing " "For DNNs to be maximally effective in aiding our understanding of splicing, "

# 3. ", they must provide information about uncertainty " -> ", provide information about uncertainty "
Replace-Text ", they must provide information about uncertainty " ", provide information about uncertainty "

# 4. ", they must be able to learn from new experiments, and they must provide insight into the underlying mechanisms of splicing. "
#    -> ", learn from new experiments, and be amenable interpretability methods. "
Replace-Text ", they must be able to learn from new experiments, and they must provide insight into the underlying mechanisms of splicing. " ", learn from new experiments, and be amenable interpretability methods. "

# 5. remove "mechanistic " before "insight into the mechanisms of splicing"
Replace-Text "develop strategies for extracting mechanistic insight into the mechanisms of splicing" "develop strategies for extracting insight into the mechanisms of splicing"

# 6. "can allow information learned in one dataset to be transferred to the other datasets. "
#    -> "will allow information learned in one dataset to transfer to the other datasets. "
Replace-Text "can allow information learned in one dataset to be transferred to the other datasets. " "will allow information learned in one dataset to transfer to the other datasets. "

# 7. " data for model fine tuning to evaluate the effectiveness of " -> " data for fine tuning to evaluate the effectiveness of "
Replace-Text " data for model fine tuning to evaluate the effectiveness of " " data for fine tuning to evaluate the effectiveness of "

# 8. "the previously learned datasets" -> "the previous datasets"
Replace-Text "the previously learned datasets" "the previous datasets"

# 9. "to specifically improve performance " -> "to improve performance "
Replace-Text "to specifically improve performance " "to improve performance "

# 10. "use neural network interpretability methods to derive mechanistic insights from splicing models. I will attribution methods such as saliency analysis and "
#     -> "use interpretability methods to derive mechanistic insights from splicing models. I will use attribution methods such as saliency maps and "
Replace-Text "use neural network interpretability methods to derive mechanistic insights from splicing models. I will attribution methods such as saliency analysis and " "use interpretability methods to derive mechanistic insights from splicing models. I will use attribution methods such as saliency maps and "

# 11. "I propose from my Ph.D. work performing massively parallel assays and my undergraduate experience majoring in molecular biology and biochemistry and performing wet lab research. "
#     -> "from my Ph.D. work performing massively parallel assays and my undergraduate experience both majoring in molecular biology and biochemistry and performing wet lab research. "
Replace-Text "I propose from my Ph.D. work performing massively parallel assays and my undergraduate experience majoring in molecular biology and biochemistry and performing wet lab research. " "from my Ph.D. work performing massively parallel assays and my undergraduate experience both majoring in molecular biology and biochemistry and performing wet lab research. "

# 12. " active learning. Through this project, I hope to build skills that will help me combine deep learning and experiments to extract insight into complex mechanisms. Further, this project will allow me to learn about the field of splicing and how to apply my expertise in massively parallel assays and computational modeling to this area. "
#     -> " active learning. Further, this project will help me to explore applications of my skills in the field of splicing. Through this project, I will build skills that will help me combine deep learning and experiments to extract insight into complex biological mechanisms. "
Replace-Text " active learning. Through this project, I hope to build skills that will help me combine deep learning and experiments to extract insight into complex mechanisms. Further, this project will allow me to learn about the field of splicing and how to apply my expertise in massively parallel assays and computational modeling to this area. " " active learning. Further, this project will help me to explore applications of my skills in the field of splicing. Through this project, I will build skills that will help me combine deep learning and experiments to extract insight into complex biological mechanisms. "
